# Auto-generated: updates cryptocurrency Price (D) and Volume(1h) (E) columns
# to match the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'33.672.12"
$ws.Range("E2").Value = "  +8.55%  "

$ws.Range("D3").Value = "'1.767.96"
$ws.Range("E3").Value = "  +4.46%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'224.58"
$ws.Range("E5").Value = "  +1.67%  "

$ws.Range("D6").Value = "'0.551"
$ws.Range("E6").Value = "  +3.22%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'30.31"
$ws.Range("E8").Value = "  +3.03%  "

$ws.Range("D9").Value = "'46.52"
$ws.Range("E9").Value = "  +3.87%  "

$ws.Range("D10").Value = "'0.277"
$ws.Range("E10").Value = "  +3.50%  "

$ws.Range("D11").Value = "'0.0657"
$ws.Range("E11").Value = "  +2.47%  "

$ws.Range("D12").Value = "'0.0923"
$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("D13").Value = "'2.020.91"
$ws.Range("E13").Value = "  +4.43%  "

$ws.Range("D14").Value = "'1.766.36"
$ws.Range("E14").Value = "  +4.25%  "

$ws.Range("D15").Value = "'0.625"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").Value = "'33.660.50"
$ws.Range("E16").Value = "  +8.45%  "

$ws.Range("D17").Value = "'9.94"
$ws.Range("E17").Value = "  -2.77%  "

$ws.Range("D18").Value = "'4.17"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("D19").Value = "'68.21"
$ws.Range("E19").Value = "  +1.95%  "

$ws.Range("D20").Value = "'250.59"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").Value = "'0.0₃0735"
$ws.Range("E21").Value = "  +2.00%  "

$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").Value = "'10.20"
$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("D24").Value = "'4.17"
$ws.Range("E24").Value = "  -2.85%  "

$ws.Range("E25").Value = "  -1.26%  "

$ws.Range("D26").Value = "'158.33"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'16.41"
$ws.Range("E27").Value = "  +2.99%  "

$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").Value = "'6.89"
$ws.Range("E29").Value = "  +2.41%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").Value = "'3.78"
$ws.Range("E31").Value = "  +6.26%  "

$ws.Range("D32").Value = "'0.0510"
$ws.Range("E32").Value = "  +1.71%  "

$ws.Range("D33").Value = "'1.19"
$ws.Range("E33").Value = "  +2.97%  "

$ws.Range("D34").Value = "'3.53"
$ws.Range("E34").Value = "  +5.20%  "

$ws.Range("D35").Value = "'1.474.94"
$ws.Range("E35").Value = "  -2.69%  "

$ws.Range("D36").Value = "'1.78"
$ws.Range("E36").Value = "  +2.69%  "

$ws.Range("E37").Value = "  +2.84%  "

$ws.Range("D38").Value = "'0.632"
$ws.Range("E38").Value = "  +2.67%  "

$ws.Range("D39").Value = "'82.90"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "'0.0184"
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("E41").Value = "  +2.37%  "

$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").Value = "'0.881"
$ws.Range("E43").Value = "  +3.96%  "

$ws.Range("D44").Value = "'2.06"
$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("D45").Value = "'0.0508"
$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("E46").Value = "  +3.15%  "

$ws.Range("D47").Value = "'1.919.36"
$ws.Range("E47").Value = "  +5.05%  "

$ws.Range("E48").Value = "  +2.93%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").Value = "'11.77"
$ws.Range("E50").Value = "  +14.27%  "

$ws.Range("D51").Value = "'50.38"
$ws.Range("E51").Value = "  -2.62%  "
